$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GA05MOAS-GL00#_R0000#_ingest_v#")

# Update the data_source column from "recovered" to "recovered_host"
$ws.Range("D2").Value = "recovered_host"
$ws.Range("D3").Value = "recovered_host"
$ws.Range("D4").Value = "recovered_host"
$ws.Range("D5").Value = "recovered_host"

# Move the active selection from B17 to C16 as in the original edit
$ws.Range("C16").Select()
